$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.754.44"
$ws.Range("E2").Value = "  +2.96%  "
$ws.Range("D3").Value = "3.135.15"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.129.86"
$ws.Range("E8").Value = "  +1.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +14.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.47%  "
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").Value = "3.651.43"
$ws.Range("E16").Value = "  +1.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "63.624.95"
$ws.Range("E18").Value = "  +2.94%  "
$ws.Range("D19").Value = "3.126.70"
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "466.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.734"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.39%  "
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("E29").Value = "  -2.00%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("D34").Value = "0.0₃0883"
$ws.Range("E34").Value = "  +11.18%  "
$ws.Range("B35").Value = "Stacks"
$ws.Range("C35").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.03%  "
$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +14.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "51.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "453.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("D43").Value = "2.904.24"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("E44").Value = "  +2.34%  "
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("E46").Value = "  +2.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.112"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.70%  "
